$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# D-column values are stored as text even though many look numeric, so we
# prefix with an apostrophe to force text entry (as a user typing into a
# General-formatted cell would need to) and then reset the style back to
# "Normal" so no stray number-format / quote-prefix styling sticks around.

$d = $ws.Range("D2")
$d.Value = "'57.752.73"
$d.Style = "Normal"
$ws.Range("E2").Value = "  -1.27%  "

$d = $ws.Range("D3")
$d.Value = "'2.442.58"
$d.Style = "Normal"
$ws.Range("E3").Value = "  -3.33%  "

$d = $ws.Range("D5")
$d.Value = "'521.54"
$d.Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "

$d = $ws.Range("D6")
$d.Value = "'130.36"
$d.Style = "Normal"
$ws.Range("E6").Value = "  -2.48%  "

$ws.Range("E7").Value = "  +0.40%  "

$ws.Range("E8").Value = "  +0.33%  "

$d = $ws.Range("D9")
$d.Value = "'2.443.33"
$d.Style = "Normal"
$ws.Range("E9").Value = "  -3.32%  "

$d = $ws.Range("D10")
$d.Value = "'0.0978"
$d.Style = "Normal"
$ws.Range("E10").Value = "  -0.12%  "

$ws.Range("E11").Value = "  -2.08%  "

$ws.Range("E12").Value = "  -4.73%  "

$d = $ws.Range("D13")
$d.Value = "'0.323"
$d.Style = "Normal"
$ws.Range("E13").Value = "  -2.95%  "

$d = $ws.Range("D14")
$d.Value = "'2.875.57"
$d.Style = "Normal"
$ws.Range("E14").Value = "  -3.31%  "

$d = $ws.Range("D15")
$d.Value = "'57.679.17"
$d.Style = "Normal"
$ws.Range("E15").Value = "  -1.39%  "

$d = $ws.Range("D16")
$d.Value = "'21.65"
$d.Style = "Normal"
$ws.Range("E16").Value = "  -2.37%  "

$ws.Range("E17").Value = "  -2.02%  "

$d = $ws.Range("D18")
$d.Value = "'2.441.07"
$d.Style = "Normal"
$ws.Range("E18").Value = "  -3.32%  "

$d = $ws.Range("D19")
$d.Value = "'10.24"
$d.Style = "Normal"
$ws.Range("E19").Value = "  -4.25%  "

$d = $ws.Range("D20")
$d.Value = "'4.12"
$d.Style = "Normal"
$ws.Range("E20").Value = "  -1.06%  "

$d = $ws.Range("D21")
$d.Value = "'316.08"
$d.Style = "Normal"
$ws.Range("E21").Value = "  -1.92%  "

$ws.Range("E22").Value = "  -1.03%  "

$ws.Range("E23").Value = "  -0.13%  "

$d = $ws.Range("D24")
$d.Value = "'64.71"
$d.Style = "Normal"
$ws.Range("E24").Value = "  -0.06%  "

$d = $ws.Range("D25")
$d.Value = "'0.400"
$d.Style = "Normal"
$ws.Range("E25").Value = "  -2.14%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("E27").Value = "  -2.32%  "

$d = $ws.Range("D28")
$d.Value = "'7.18"
$d.Style = "Normal"
$ws.Range("E28").Value = "  -3.08%  "

$d = $ws.Range("D29")
$d.Value = "'174.59"
$d.Style = "Normal"
$ws.Range("E29").Value = "  +3.94%  "

$d = $ws.Range("D30")
$d.Value = "'0.0₃0733"
$d.Style = "Normal"
$ws.Range("E30").Value = "  -3.17%  "

$d = $ws.Range("D31")
$d.Value = "'1.69"
$d.Style = "Normal"
$ws.Range("E31").Value = "  -2.37%  "

$ws.Range("E33").Value = "  -5.44%  "

$d = $ws.Range("D34")
$d.Value = "'0.999"
$d.Style = "Normal"
$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("E35").Value = "  -0.03%  "

$d = $ws.Range("D36")
$d.Value = "'17.78"
$d.Style = "Normal"
$ws.Range("E36").Value = "  -2.05%  "

$ws.Range("E37").Value = "  -6.27%  "

$d = $ws.Range("D38")
$d.Value = "'3.76"
$d.Style = "Normal"
$ws.Range("E38").Value = "  -4.87%  "

$d = $ws.Range("D39")
$d.Value = "'36.11"
$d.Style = "Normal"
$ws.Range("E39").Value = "  -0.87%  "

$d = $ws.Range("D40")
$d.Value = "'1.44"
$d.Style = "Normal"
$ws.Range("E40").Value = "  -2.61%  "

$d = $ws.Range("D41")
$d.Value = "'0.791"
$d.Style = "Normal"
$ws.Range("E41").Value = "  +1.90%  "

$d = $ws.Range("D42")
$d.Value = "'3.41"
$d.Style = "Normal"
$ws.Range("E42").Value = "  -2.64%  "

$d = $ws.Range("D43")
$d.Value = "'261.97"
$d.Style = "Normal"
$ws.Range("E43").Value = "  -5.78%  "

$d = $ws.Range("D44")
$d.Value = "'0.584"
$d.Style = "Normal"
$ws.Range("E44").Value = "  -2.69%  "

$ws.Range("E45").Value = "  -4.21%  "

$d = $ws.Range("D46")
$d.Value = "'0.0920"
$d.Style = "Normal"
$ws.Range("E46").Value = "  -0.33%  "

$d = $ws.Range("D47")
$d.Value = "'122.13"
$d.Style = "Normal"
$ws.Range("E47").Value = "  -6.35%  "

$d = $ws.Range("D48")
$d.Value = "'0.0494"
$d.Style = "Normal"
$ws.Range("E48").Value = "  -1.90%  "

$d = $ws.Range("D49")
$d.Value = "'0.0210"
$d.Style = "Normal"
$ws.Range("E49").Value = "  -1.74%  "

$d = $ws.Range("D50")
$d.Value = "'16.94"
$d.Style = "Normal"
$ws.Range("E50").Value = "  -4.90%  "

$d = $ws.Range("D51")
$d.Value = "'16.22"
$d.Style = "Normal"
$ws.Range("E51").Value = "  -4.27%  "
